$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (file name) for all data rows (2-17) to the new raw file name,
# and add the new "cell_equivalents" value (16) in column C for the same rows.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = "CellenONE_I3T_NEM_SC_Chip1_C1.raw"
    $ws.Cells.Item($r, 3).Value = 16
}

# Update the active selection to A10
$ws.Range("A10").Select()
